# "further cleaning to metadata"
#  - fix the typo'd libraryProtocol code E7760 -> E7420 (all K2:K27 cells share it)
#  - give the libraryProtocol column (K2:K27) its own dedicated font/style
#  - turn the literal FALSE booleans in L2:L27 (roboticLibraryPrep) into real =FALSE() formulas
#  - follow the selection that was left on the sheet (K2:K27 instead of L2:L27)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the libraryProtocol text used by every row.
$ws.Cells.Replace("E7760", "E7420") | Out-Null

# 2) Give K2:K27 (libraryProtocol) its own font (same family, larger size,
#    explicit black) so it gets a dedicated style, matching the new
#    cellXf/font pair added to styles.xml.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 11).Font.Size = 11
    $ws.Cells.Item($r, 11).Font.Color = 0
}

# 3) Replace the literal FALSE booleans in L2:L27 with real formulas.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 12).Formula = "=FALSE()"
}

# 4) Match the saved selection (K2:K27 is now the interesting column).
$ws.Range("K2:K27").Select() | Out-Null
